$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of column C (codeforiati:group-code) and column D
# (codeforiati:group-name) for every used row, so that column C now
# holds the group-name values and column D holds the group-code values
# (including the header row itself).

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
